$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Technology_Group values for column B (row 1 = header "Technology_Group",
# rows 2-90 = the per-technology group assignment)
$groups = @(
    'Technology_Group',
    'Mobile Motive Power',
    'Heating/Cooling',
    'Heating/Cooling',
    'Heating/Cooling',
    'Stationary Motive Power',
    'Stationary Motive Power',
    'Electronics and Lighting',
    'Heating/Cooling',
    'Mobile Motive Power',
    'Stationary Motive Power',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Stationary Motive Power',
    'Stationary Motive Power',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Stationary Motive Power',
    'Heating/Cooling',
    'Heating/Cooling',
    'Heating/Cooling',
    'Other',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Heating/Cooling',
    'Electronics and Lighting',
    'Heating/Cooling',
    'Electronics and Lighting',
    'Electronics and Lighting',
    'Mobile Motive Power',
    'Heating/Cooling',
    'Heating/Cooling',
    'Electronics and Lighting',
    'Heating/Cooling',
    'Fuel Production',
    'Heating/Cooling',
    'Mobile Motive Power',
    'Stationary Motive Power',
    'Electronics and Lighting',
    'Electronics and Lighting',
    'Heating/Cooling',
    'Stationary Motive Power',
    'Feedstock',
    'Electronics and Lighting',
    'Heating/Cooling',
    'Heating/Cooling',
    'Stationary Motive Power',
    'Storage',
    'Storage',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production',
    'Heating/Cooling',
    'Heating/Cooling',
    'Heating/Cooling',
    'Heating/Cooling',
    'Heating/Cooling',
    'Heating/Cooling',
    'Heating/Cooling',
    'Heating/Cooling',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Mobile Motive Power',
    'Fuel Production',
    'Fuel Production',
    'Fuel Production'
)

for ($i = 0; $i -lt $groups.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $groups[$i]
}

# Widen column B so the new, longer group names fit
$ws.Columns.Item(2).ColumnWidth = 22.5

# Move the active selection to B4 (matches the saved view state)
$ws.Range("B4").Select()
